$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'https://www.google.com/url?url=https://www.lebiscuit.com.br/iphone-12-apple-64gb-tela-de-6-1-polegadas-camera-12mp-ios-16597al08l832478/p%3Fidsku%3D2147376583&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwjK69Gf5tWCAxXeALkGHSUvAroQgOUECP8N&usg=AOvVaw1sxxN35YJdbbmagNBLXr-S'

$ws.Range("B3").Value = 3379
$ws.Range("C3").Value = 'https://www.google.com/url?url=https://www.trocafone.com.br/iphone-12-64gb-azul---excelente-70613/p%3Fidsku%3D39869&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwjK69Gf5tWCAxXeALkGHSUvAroQguUECKoP&usg=AOvVaw0YyKGqD3oITCUweNcmrHnF'

$ws.Range("A4").Value = 'usado: iphone 12 64gb vermelho excelente - trocafone'
$ws.Range("B4").Value = 3209
$ws.Range("C4").Value = 'https://www.google.com/url?url=%23&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwjK69Gf5tWCAxXeALkGHSUvAroQgOUECOwP&usg=AOvVaw0ELP03YRLky63eVe8qMhUy'

$ws.Range("A5").Value = 'usado: iphone 12 64gb branco excelente - trocafone'
$ws.Range("B5").Value = 3209
$ws.Range("C5").Value = 'https://www.google.com/url?url=https://www.casasbahia.com.br/usado-iphone-12-64gb-branco-excelente-trocafone-1561610772/p/1561610772%3Futm_medium%3DCpc%26utm_source%3Dgoogle_freelisting%26IdSku%3D1561610772%26idLojista%3D33013%26tipoLojista%3D3P&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwjK69Gf5tWCAxXeALkGHSUvAroQgOUECOoQ&usg=AOvVaw2JaTFQyj-bFxPV3GFIdaaY'

$ws.Range("C6").Value = 'https://www.google.com/url?url=https://www.carrefour.com.br/vitrine-iphone-12-preto-64gb-mp929837648/p&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwjK69Gf5tWCAxXeALkGHSUvAroQgOUECLoR&usg=AOvVaw2kPrazmE3tzY-rx2dyp3FC'

$ws.Range("A7").Value = 'iphone 12 64gb | celular apple | usado'
$ws.Range("B7").Value = 3315
$ws.Range("C7").Value = 'https://www.google.com/url?url=https://www.enjoei.com.br/p/iphone-12-64gb-93394889%3Fg_campaign%3Dgoogle_shopping%26srsltid%3DAfmBOortXtLZvDOqq1sItyjrGQONVn_Qk310LKRAjzTsuSMcG2VLECaWKEo&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwjK69Gf5tWCAxXeALkGHSUvAroQgOUECMoR&usg=AOvVaw2UmlbIqLpkQKm1KLx94GKM'

$ws.Range("C9").Value = 'https://www.google.com/url?url=https://gamernareal.com.br/product/placa-de-video-msi-geforce-rtx-2060-ventus-gp/%3Futm_source%3DGoogle%2520Shopping%26utm_campaign%3DGoogle%2520Menor%2520Preco%26utm_medium%3Dcpc%26utm_term%3D7362&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwin_Mum5tWCAxX3HrkGHYMVCmkQguUECJkL&usg=AOvVaw3OOqLcTtPly5vSBRB_ItCU'

$ws.Range("C10").Value = 'https://www.google.com/url?url=https://www.casasbahia.com.br/placa-video-pcwinmax-nvidia-geforce-rtx-2060-12gb-192bits-1559067608/p/1559067608%3Futm_medium%3DCpc%26utm_source%3Dgoogle_freelisting%26IdSku%3D1559067608%26idLojista%3D12231%26tipoLojista%3D3P&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwin_Mum5tWCAxX3HrkGHYMVCmkQguUECOwM&usg=AOvVaw0eEqjWpPnNQi9mQBYTDd3P'

$ws.Range("A11").Value = 'placa de vídeo pcyes geforce rtx 2060 6gb gddr6 192 bits'
$ws.Range("B11").Value = 2199
$ws.Range("C11").Value = 'https://www.google.com/url?url=https://www.fujioka.com.br/placa-de-video-pcyes-nvidia-geforce-rtx-2060-6gb-gddr6-192-bits-dual-fan-pax2060rtx6g192bdf-801221/p%3Fidsku%3D27526&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwin_Mum5tWCAxX3HrkGHYMVCmkQguUECPwN&usg=AOvVaw1U7jE4_D3InWomTlfOLbds'

$ws.Range("A12").Value = 'nova placa de vídeo asus dual geforce rtx 2060 evo oc edition 12gddr6 gb memória'
$ws.Range("B12").Value = 2028.16
$ws.Range("C12").Value = 'https://www.google.com/url?url=https://www.ebay.com/itm/255821325545%3Fchn%3Dps%26mkevt%3D1%26mkcid%3D28&rct=j&q=&esrc=s&opi=95576897&sa=U&ved=0ahUKEwin_Mum5tWCAxX3HrkGHYMVCmkQgOUECMwO&usg=AOvVaw1cO4NyHqVsERmt2qAF0R2H'

$ws.Range("A13").Value = 'evga placa de vídeo nvidia geforce rtx 2060 6gb gddr6'
$ws.Range("B13").Value = 1990
$ws.Range("C13").Value = 'https://www.buscape.com.br/lead?oid=1139015123&channel=86&index=2&searchterm=rtx%202060'

$ws.Range("A14").Value = 'evga placa gráfica geforce rtx 2060 sc gaming, 6gb gddr6, hdb fan 06g-p4-2062-kr'
$ws.Range("B14").Value = 1979.9
$ws.Range("C14").Value = 'https://www.buscape.com.br/lead?oid=841411028&channel=86&index=3&searchterm=rtx%202060'

$ws.Range("A15:C15").EntireRow.Delete()
